$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 201-202, pushing existing rows 201+ down by 2
$ws.Rows("201:202").Insert()

# New row 201 data
$ws.Cells.Item(201, 1).Value = 10
$ws.Cells.Item(201, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(201, 3).Value = "La Araucanía"
$ws.Cells.Item(201, 4).Value = 44841
$ws.Cells.Item(201, 4).NumberFormat = $ws.Cells.Item(200, 4).NumberFormat
$ws.Cells.Item(201, 5).Value = 9
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100101
$ws.Cells.Item(201, 8).Value = "Berries"
$ws.Cells.Item(201, 9).Value = 100112025
$ws.Cells.Item(201, 10).Value = "Frutilla"
$ws.Cells.Item(201, 11).Value = "Sin especificar"
$ws.Cells.Item(201, 12).Value = "Primera"
$ws.Cells.Item(201, 13).Value = 215
$ws.Cells.Item(201, 14).Value = 14000
$ws.Cells.Item(201, 15).Value = 16000
$ws.Cells.Item(201, 16).Value = 15163
$ws.Cells.Item(201, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(201, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(201, 19).Value = 2166
$ws.Cells.Item(201, 20).Value = 7

# New row 202 data
$ws.Cells.Item(202, 1).Value = 10
$ws.Cells.Item(202, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(202, 3).Value = "La Araucanía"
$ws.Cells.Item(202, 4).Value = 44841
$ws.Cells.Item(202, 4).NumberFormat = $ws.Cells.Item(200, 4).NumberFormat
$ws.Cells.Item(202, 5).Value = 9
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100101
$ws.Cells.Item(202, 8).Value = "Berries"
$ws.Cells.Item(202, 9).Value = 100112025
$ws.Cells.Item(202, 10).Value = "Frutilla"
$ws.Cells.Item(202, 11).Value = "Sin especificar"
$ws.Cells.Item(202, 12).Value = "Segunda"
$ws.Cells.Item(202, 13).Value = 155
$ws.Cells.Item(202, 14).Value = 8000
$ws.Cells.Item(202, 15).Value = 8000
$ws.Cells.Item(202, 16).Value = 8000
$ws.Cells.Item(202, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(202, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(202, 19).Value = 1143
$ws.Cells.Item(202, 20).Value = 7
